# Fragenbaum.xlsx update:
#  - Add "Ja" answers for rows 26, 31 and 33 (column D)
#  - Move the sheet's scroll/selection so D34 (just past the last data
#    row) is the active cell, with the view scrolled down near row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "Ja"
$ws.Range("D31").Value = "Ja"
$ws.Range("D33").Value = "Ja"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D34").Select()
